# Add cost/price default rows for each country, plus a bold header row.
# "added all costs - it's overly busy but works"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make the header row bold (introduces a new bold font + cellXf).
$ws.Range("A1:C1").Font.Bold = $true

# country, variable, value, useTwoDecimalFormat
$data = @(
    @("China",   "price_egg",     0.53,  $true),
    @("China",   "price_spent",   0.21,  $true),
    @("China",   "price_manure",  0.11,  $true),
    @("Vietnam", "price_egg",     1.2,   $true),
    @("Vietnam", "price_spent",   0.09,  $true),
    @("Vietnam", "price_manure",  0.2,   $true),
    @("China",   "cost_feed",     1000,  $false),
    @("China",   "cost_labor",    1.2,   $false),
    @("China",   "cost_pullet",   1.1,   $false),
    @("China",   "cost_equip",    3,     $false),
    @("China",   "cost_litter",   2,     $false),
    @("China",   "cost_vet",      3,     $false),
    @("Vietnam", "cost_feed",     800,   $false),
    @("Vietnam", "cost_labor",    1.8,   $false),
    @("Vietnam", "cost_pullet",   3,     $false),
    @("Vietnam", "cost_equip",    0.4,   $false),
    @("Vietnam", "cost_litter",   1.1,   $false),
    @("Vietnam", "cost_vet",      0.9,   $false),
    @("Vietnam", "cost_land",     40000, $false),
    @("Vietnam", "cost_office",   50000, $false),
    @("China",   "cost_land",     30000, $false),
    @("China",   "cost_office",   34000, $false)
)

$row = 14
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $cell = $ws.Cells.Item($row, 3)
    $cell.Value = $item[2]
    if ($item[3]) {
        $cell.NumberFormat = "0.00"
    }
    $row = $row + 1
}

# Mirror the saved selection/view state from the edit (last cell below
# the new data).
$ws.Range("C36").Select()
